$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 58906
$ws.Cells.Item(17, 10).Value = 58906
$ws.Cells.Item(17, 12).Value = 176718
$ws.Cells.Item(17, 14).Value = -177054
$ws.Cells.Item(33, 8).Value = 543.02856
$ws.Cells.Item(33, 9).Value = 461.31818
$ws.Cells.Item(33, 10).Value = 681.3077
$ws.Cells.Item(33, 11).Value = 461.31818
$ws.Cells.Item(33, 12).Value = 681.3077
$ws.Cells.Item(33, 13).Value = -232.31818
$ws.Cells.Item(33, 14).Value = -1139.3077
$ws.Cells.Item(53, 8).Value = 416.81482
$ws.Cells.Item(53, 9).Value = 572.125
$ws.Cells.Item(53, 10).Value = 351.42105
$ws.Cells.Item(53, 11).Value = 572.125
$ws.Cells.Item(53, 12).Value = 351.42105
$ws.Cells.Item(53, 13).Value = 64.875
$ws.Cells.Item(53, 14).Value = -1625.42105
$ws.Cells.Item(62, 8).Value = 5105
$ws.Cells.Item(62, 10).Value = 5105
$ws.Cells.Item(62, 12).Value = 5105
$ws.Cells.Item(62, 14).Value = -6353
$ws.Cells.Item(65, 8).Value = 5105
$ws.Cells.Item(65, 10).Value = 5105
$ws.Cells.Item(65, 12).Value = 25525
$ws.Cells.Item(65, 14).Value = -31765
$ws.Cells.Item(98, 8).Value = 4256.7646
$ws.Cells.Item(98, 9).Value = 4655.4165
$ws.Cells.Item(98, 10).Value = 3300
$ws.Cells.Item(98, 11).Value = 4655.4165
$ws.Cells.Item(98, 12).Value = 3300
$ws.Cells.Item(98, 13).Value = -3157.4165
$ws.Cells.Item(98, 14).Value = -6296
$ws.Cells.Item(105, 8).Value = 29975
$ws.Cells.Item(105, 10).Value = 29975
$ws.Cells.Item(105, 12).Value = 29975
$ws.Cells.Item(105, 14).Value = -36963
$ws.Cells.Item(122, 8).Value = 4256.7646
$ws.Cells.Item(122, 9).Value = 4655.4165
$ws.Cells.Item(122, 10).Value = 3300
$ws.Cells.Item(122, 11).Value = 13966.2495
$ws.Cells.Item(122, 12).Value = 9900
$ws.Cells.Item(122, 13).Value = -11516.2495
$ws.Cells.Item(122, 14).Value = -14800
$ws.Cells.Item(132, 8).Value = 5717730.5
$ws.Cells.Item(132, 9).Value = 7410786
$ws.Cells.Item(132, 11).Value = 22232358
$ws.Cells.Item(132, 13).Value = -22229828
$ws.Cells.Item(137, 8).Value = 3666.2708
$ws.Cells.Item(137, 9).Value = 4152.207
$ws.Cells.Item(137, 10).Value = 2924.5789
$ws.Cells.Item(137, 11).Value = 12456.621
$ws.Cells.Item(137, 12).Value = 8773.736699999999
$ws.Cells.Item(137, 13).Value = -9906.621000000001
$ws.Cells.Item(137, 14).Value = -13873.7367
$ws.Cells.Item(138, 8).Value = 2267.1865
$ws.Cells.Item(138, 9).Value = 1460.0968
$ws.Cells.Item(138, 10).Value = 3160.75
$ws.Cells.Item(138, 11).Value = 4380.2904
$ws.Cells.Item(138, 12).Value = 9482.25
$ws.Cells.Item(138, 13).Value = 759.7096000000001
$ws.Cells.Item(138, 14).Value = -19762.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 5683306.5
$ws.Cells.Item(2, 9).Value = 9616661
$ws.Cells.Item(2, 10).Value = 1795.4445
$ws.Cells.Item(2, 11).Value = 9616661
$ws.Cells.Item(2, 12).Value = 1795.4445
$ws.Cells.Item(2, 13).Value = -9616548
$ws.Cells.Item(2, 14).Value = -2021.4445
$ws.Cells.Item(32, 8).Value = 6386.46
$ws.Cells.Item(32, 9).Value = 5655.125
$ws.Cells.Item(32, 10).Value = 11749.583
$ws.Cells.Item(32, 11).Value = 5655.125
$ws.Cells.Item(32, 12).Value = 11749.583
$ws.Cells.Item(32, 13).Value = -5368.125
$ws.Cells.Item(32, 14).Value = -12323.583
$ws.Cells.Item(45, 8).Value = 1452.4
$ws.Cells.Item(45, 9).Value = 1067.2069
$ws.Cells.Item(45, 10).Value = 3314.1667
$ws.Cells.Item(45, 11).Value = 1067.2069
$ws.Cells.Item(45, 12).Value = 3314.1667
$ws.Cells.Item(45, 13).Value = -690.2068999999999
$ws.Cells.Item(45, 14).Value = -4068.1667
$ws.Cells.Item(46, 8).Value = 5753.4
$ws.Cells.Item(46, 9).Value = 0
$ws.Cells.Item(46, 10).Value = 5753.4
$ws.Cells.Item(46, 11).Value = 0
$ws.Cells.Item(46, 12).Value = 5753.4
$ws.Cells.Item(46, 13).ClearContents()
$ws.Cells.Item(46, 14).Value = -6391.4
$ws.Cells.Item(116, 8).Value = 5683306.5
$ws.Cells.Item(116, 9).Value = 9616661
$ws.Cells.Item(116, 10).Value = 1795.4445
$ws.Cells.Item(116, 11).Value = 9616661
$ws.Cells.Item(116, 12).Value = 1795.4445
$ws.Cells.Item(116, 13).Value = -9614367
$ws.Cells.Item(116, 14).Value = -6383.4445
$ws.Cells.Item(122, 8).Value = 1915.9412
$ws.Cells.Item(122, 9).Value = 1412
$ws.Cells.Item(122, 10).Value = 2554.2666
$ws.Cells.Item(122, 11).Value = 4236
$ws.Cells.Item(122, 12).Value = 7662.7998
$ws.Cells.Item(122, 13).Value = -1786
$ws.Cells.Item(122, 14).Value = -12562.7998
$ws.Cells.Item(132, 8).Value = 1792.2858
$ws.Cells.Item(132, 9).Value = 1473.0344
$ws.Cells.Item(132, 11).Value = 4419.1032
$ws.Cells.Item(132, 13).Value = -1889.1032

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 5683306.5
$ws.Cells.Item(3, 9).Value = 9616661
$ws.Cells.Item(3, 10).Value = 1795.4445
$ws.Cells.Item(3, 11).Value = 9616661
$ws.Cells.Item(3, 12).Value = 1795.4445
$ws.Cells.Item(3, 13).Value = -9616547
$ws.Cells.Item(3, 14).Value = -2023.4445
$ws.Cells.Item(134, 8).Value = 1936.7693
$ws.Cells.Item(134, 9).Value = 1264.8334
$ws.Cells.Item(134, 10).Value = 10000
$ws.Cells.Item(134, 11).Value = 3794.5002
$ws.Cells.Item(134, 12).Value = 30000
$ws.Cells.Item(134, 13).Value = -1259.5002
$ws.Cells.Item(134, 14).Value = -35070

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 2739.4
$ws.Cells.Item(16, 9).Value = 1923.75
$ws.Cells.Item(16, 10).Value = 3283.1667
$ws.Cells.Item(16, 11).Value = 1923.75
$ws.Cells.Item(16, 12).Value = 3283.1667
$ws.Cells.Item(16, 13).Value = -1636.75
$ws.Cells.Item(16, 14).Value = -3857.1667
$ws.Cells.Item(58, 8).Value = 13517533
$ws.Cells.Item(58, 9).Value = 1697.591
$ws.Cells.Item(58, 10).Value = 33340760
$ws.Cells.Item(58, 11).Value = 1697.591
$ws.Cells.Item(58, 12).Value = 33340760
$ws.Cells.Item(58, 13).Value = -1494.591
$ws.Cells.Item(58, 14).Value = -33341166
$ws.Cells.Item(105, 8).Value = 1905.625
$ws.Cells.Item(105, 9).Value = 1499.1666
$ws.Cells.Item(105, 11).Value = 1499.1666
$ws.Cells.Item(105, 13).Value = 247.8334
$ws.Cells.Item(113, 8).Value = 2739.4
$ws.Cells.Item(113, 9).Value = 1923.75
$ws.Cells.Item(113, 10).Value = 3283.1667
$ws.Cells.Item(113, 11).Value = 1923.75
$ws.Cells.Item(113, 12).Value = 3283.1667
$ws.Cells.Item(113, 13).Value = 246.25
$ws.Cells.Item(113, 14).Value = -7623.1667
$ws.Cells.Item(132, 8).Value = 2238.9744
$ws.Cells.Item(132, 9).Value = 2105.4546
$ws.Cells.Item(132, 10).Value = 2411.7646
$ws.Cells.Item(132, 11).Value = 6316.3638
$ws.Cells.Item(132, 12).Value = 7235.293799999999
$ws.Cells.Item(132, 13).Value = -3786.3638
$ws.Cells.Item(132, 14).Value = -12295.2938
$ws.Cells.Item(134, 8).Value = 4154.136
$ws.Cells.Item(134, 9).Value = 4512.7334
$ws.Cells.Item(134, 10).Value = 3385.7144
$ws.Cells.Item(134, 11).Value = 13538.2002
$ws.Cells.Item(134, 12).Value = 10157.1432
$ws.Cells.Item(134, 13).Value = -11003.2002
$ws.Cells.Item(134, 14).Value = -15227.1432
$ws.Cells.Item(136, 8).Value = 13517533
$ws.Cells.Item(136, 9).Value = 1697.591
$ws.Cells.Item(136, 10).Value = 33340760
$ws.Cells.Item(136, 11).Value = 5092.772999999999
$ws.Cells.Item(136, 12).Value = 100022280
$ws.Cells.Item(136, 13).Value = -2542.772999999999
$ws.Cells.Item(136, 14).Value = -100027380

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(86, 8).Value = 1775.75
$ws.Cells.Item(86, 9).Value = 1500
$ws.Cells.Item(86, 10).Value = 1867.6666
$ws.Cells.Item(86, 11).Value = 4500
$ws.Cells.Item(86, 12).Value = 5602.9998
$ws.Cells.Item(86, 13).Value = -3314
$ws.Cells.Item(86, 14).Value = -7974.9998
$ws.Cells.Item(87, 8).Value = 12555.444
$ws.Cells.Item(87, 9).Value = 1999.5
$ws.Cells.Item(87, 10).Value = 15571.429
$ws.Cells.Item(87, 11).Value = 5998.5
$ws.Cells.Item(87, 12).Value = 46714.287
$ws.Cells.Item(87, 13).Value = -4750.5
$ws.Cells.Item(87, 14).Value = -49210.287
$ws.Cells.Item(88, 8).Value = 4114.6665
$ws.Cells.Item(88, 10).Value = 4114.6665
$ws.Cells.Item(88, 12).Value = 12343.9995
$ws.Cells.Item(88, 14).Value = -13199.9995
$ws.Cells.Item(89, 8).Value = 1775.75
$ws.Cells.Item(89, 9).Value = 1500
$ws.Cells.Item(89, 10).Value = 1867.6666
$ws.Cells.Item(89, 11).Value = 13500
$ws.Cells.Item(89, 12).Value = 16808.9994
$ws.Cells.Item(89, 13).Value = -7572
$ws.Cells.Item(89, 14).Value = -28664.9994
$ws.Cells.Item(90, 8).Value = 12555.444
$ws.Cells.Item(90, 9).Value = 1999.5
$ws.Cells.Item(90, 10).Value = 15571.429
$ws.Cells.Item(90, 11).Value = 17995.5
$ws.Cells.Item(90, 12).Value = 140142.861
$ws.Cells.Item(90, 13).Value = -11755.5
$ws.Cells.Item(90, 14).Value = -152622.861
$ws.Cells.Item(91, 8).Value = 4114.6665
$ws.Cells.Item(91, 10).Value = 4114.6665
$ws.Cells.Item(91, 12).Value = 12343.9995
$ws.Cells.Item(91, 14).Value = -15307.9995
$ws.Cells.Item(92, 8).Value = 840.2
$ws.Cells.Item(92, 9).Value = 0
$ws.Cells.Item(92, 10).Value = 840.2
$ws.Cells.Item(92, 11).Value = 0
$ws.Cells.Item(92, 12).Value = 2520.6
$ws.Cells.Item(92, 13).ClearContents()
$ws.Cells.Item(92, 14).Value = -5016.6
$ws.Cells.Item(98, 8).Value = 250.94444
$ws.Cells.Item(98, 9).Value = 180.75
$ws.Cells.Item(98, 10).Value = 391.33334
$ws.Cells.Item(98, 11).Value = 542.25
$ws.Cells.Item(98, 12).Value = 1174.00002
$ws.Cells.Item(98, 13).Value = 955.75
$ws.Cells.Item(98, 14).Value = -4170.000019999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 34785.902
$ws.Cells.Item(102, 9).Value = 1785.55
$ws.Cells.Item(102, 10).Value = 94786.55
$ws.Cells.Item(102, 11).Value = 1785.55
$ws.Cells.Item(102, 12).Value = 94786.55
$ws.Cells.Item(102, 13).Value = -163.55
$ws.Cells.Item(102, 14).Value = -98030.55
$ws.Cells.Item(121, 8).Value = 24000
$ws.Cells.Item(121, 10).Value = 24000
$ws.Cells.Item(121, 12).Value = 24000
$ws.Cells.Item(121, 13).Value = -27494
$ws.Cells.Item(132, 8).Value = 2830.1052
$ws.Cells.Item(132, 9).Value = 2591.0881
$ws.Cells.Item(132, 10).Value = 3183.4348
$ws.Cells.Item(132, 11).Value = 7773.2643
$ws.Cells.Item(132, 12).Value = 9550.304400000001
$ws.Cells.Item(132, 13).Value = -5243.2643
$ws.Cells.Item(132, 14).Value = -14610.3044

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3573355.2
$ws.Cells.Item(7, 9).Value = 7144194.5
$ws.Cells.Item(7, 10).Value = 2515.9285
$ws.Cells.Item(7, 11).Value = 7144194.5
$ws.Cells.Item(7, 12).Value = 2515.9285
$ws.Cells.Item(7, 13).Value = -7144082.5
$ws.Cells.Item(7, 14).Value = -2739.9285
$ws.Cells.Item(22, 8).Value = 200002270
$ws.Cells.Item(22, 9).Value = 250000340
$ws.Cells.Item(22, 10).Value = 10000
$ws.Cells.Item(22, 11).Value = 250000340
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = -250000045
$ws.Cells.Item(22, 14).Value = -10590
$ws.Cells.Item(27, 8).Value = 200002270
$ws.Cells.Item(27, 9).Value = 250000340
$ws.Cells.Item(27, 10).Value = 10000
$ws.Cells.Item(27, 11).Value = 250000340
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = -250000233
$ws.Cells.Item(27, 14).Value = -10214
$ws.Cells.Item(40, 8).Value = 2750
$ws.Cells.Item(40, 9).Value = 1500
$ws.Cells.Item(40, 10).Value = 4000
$ws.Cells.Item(40, 11).Value = 1500
$ws.Cells.Item(40, 12).Value = 4000
$ws.Cells.Item(40, 13).Value = -1364
$ws.Cells.Item(40, 14).Value = -4272
$ws.Cells.Item(126, 8).Value = 3573355.2
$ws.Cells.Item(126, 9).Value = 7144194.5
$ws.Cells.Item(126, 10).Value = 2515.9285
$ws.Cells.Item(126, 11).Value = 21432583.5
$ws.Cells.Item(126, 12).Value = 7547.7855
$ws.Cells.Item(126, 13).Value = -21430113.5
$ws.Cells.Item(126, 14).Value = -12487.7855
$ws.Cells.Item(132, 8).Value = 2934.6
$ws.Cells.Item(132, 9).Value = 2185.3
$ws.Cells.Item(132, 11).Value = 6555.900000000001
$ws.Cells.Item(132, 13).Value = -4025.900000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 443.44446
$ws.Cells.Item(100, 9).Value = 411.375
$ws.Cells.Item(100, 10).Value = 700
$ws.Cells.Item(100, 11).Value = 822.75
$ws.Cells.Item(100, 12).Value = 1400
$ws.Cells.Item(100, 13).Value = -281.75
$ws.Cells.Item(100, 14).Value = -2482
$ws.Cells.Item(126, 8).Value = 2082.1316
$ws.Cells.Item(126, 10).Value = 3107.0833
$ws.Cells.Item(126, 12).Value = 9321.249899999999
$ws.Cells.Item(126, 14).Value = -14261.2499
$ws.Cells.Item(132, 8).Value = 15530.718
$ws.Cells.Item(132, 9).Value = 2358.889
$ws.Cells.Item(132, 11).Value = 7076.667
$ws.Cells.Item(132, 13).Value = -4546.667
